$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("B5").Value = 6102174
$ws.Range("E5").Value = 'Breidablik'
$ws.Range("F5").Value = 'Vikingur Reykjavik'
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 'D'
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 3.6
$ws.Range("L5").Value = 3.1
$ws.Range("M5").Value = 2.25
$ws.Range("N5").Value = 3.3
$ws.Range("O5").Value = 2.8
$ws.Range("P5").Value = -0.25
$ws.Range("Q5").Value = 2.025
$ws.Range("R5").Value = 1.775
$ws.Range("S5").Value = 2.75
$ws.Range("T5").Value = 1.9
$ws.Range("U5").Value = 1.9
$ws.Range("V5").Value = -1
$ws.Range("W5").Value = 2.3
$ws.Range("X5").Value = -1
$ws.Range("Y5").Value = -0.5
$ws.Range("Z5").Value = 0.3875
$ws.Range("AA5").Value = 0.8999999999999999
$ws.Range("AB5").Value = -1

# Row 6
$ws.Range("B6").Value = 6102346
$ws.Range("E6").Value = 'Valur Reykjavik'
$ws.Range("F6").Value = 'FH Hafnarfjordur'
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 'D'
$ws.Range("J6").Value = 1.666
$ws.Range("K6").Value = 4
$ws.Range("L6").Value = 4.2
$ws.Range("M6").Value = 1.615
$ws.Range("N6").Value = 4.2
$ws.Range("O6").Value = 4.333
$ws.Range("P6").Value = -1
$ws.Range("Q6").Value = 2.05
$ws.Range("R6").Value = 1.8
$ws.Range("S6").Value = 3.25
$ws.Range("T6").Value = 2.025
$ws.Range("U6").Value = 1.825
$ws.Range("V6").Value = -1
$ws.Range("W6").Value = 3.2
$ws.Range("X6").Value = -1
$ws.Range("Y6").Value = -1
$ws.Range("Z6").Value = 0.8
$ws.Range("AA6").Value = -1
$ws.Range("AB6").Value = 0.825

# Row 7
$ws.Range("B7").Value = 6102175
$ws.Range("E7").Value = 'Fram Reykjavik'
$ws.Range("F7").Value = 'Keflavik'
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 'H'
$ws.Range("J7").Value = 1.909
$ws.Range("K7").Value = 3.75
$ws.Range("L7").Value = 3.25
$ws.Range("M7").Value = 1.909
$ws.Range("N7").Value = 3.6
$ws.Range("O7").Value = 3.5
$ws.Range("P7").Value = -0.5
$ws.Range("Q7").Value = 1.95
$ws.Range("R7").Value = 1.85
$ws.Range("S7").Value = 2.75
$ws.Range("T7").Value = 1.8
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 0.909
$ws.Range("W7").Value = -1
$ws.Range("X7").Value = -1
$ws.Range("Y7").Value = 0.95
$ws.Range("Z7").Value = -1
$ws.Range("AA7").Value = 0.8
$ws.Range("AB7").Value = -1

# Row 12
$ws.Range("B12").Value = 6102180
$ws.Range("E12").Value = 'Keflavik'
$ws.Range("F12").Value = 'Stjarnan'
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 'D'
$ws.Range("J12").Value = 3.4
$ws.Range("K12").Value = 3.6
$ws.Range("L12").Value = 1.85
$ws.Range("M12").Value = 4.2
$ws.Range("N12").Value = 3.5
$ws.Range("O12").Value = 1.727
$ws.Range("P12").Value = 0.75
$ws.Range("Q12").Value = 1.8
$ws.Range("R12").Value = 2
$ws.Range("S12").Value = 2.75
$ws.Range("T12").Value = 1.975
$ws.Range("U12").Value = 1.825
$ws.Range("V12").Value = -1
$ws.Range("W12").Value = 2.5
$ws.Range("X12").Value = -1
$ws.Range("Y12").Value = 0.8
$ws.Range("Z12").Value = -1
$ws.Range("AA12").Value = -1
$ws.Range("AB12").Value = 0.825

# Row 13
$ws.Range("B13").Value = 6102179
$ws.Range("E13").Value = 'Vikingur Reykjavik'
$ws.Range("F13").Value = 'Fram Reykjavik'
$ws.Range("G13").Value = 3
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 'H'
$ws.Range("J13").Value = 1.4
$ws.Range("K13").Value = 4.5
$ws.Range("L13").Value = 5.5
$ws.Range("M13").Value = 1.3
$ws.Range("N13").Value = 5.25
$ws.Range("O13").Value = 6.5
$ws.Range("P13").Value = -1.5
$ws.Range("Q13").Value = 1.85
$ws.Range("R13").Value = 1.95
$ws.Range("S13").Value = 3.5
$ws.Range("T13").Value = 1.775
$ws.Range("U13").Value = 1.925
$ws.Range("V13").Value = 0.3
$ws.Range("W13").Value = -1
$ws.Range("X13").Value = -1
$ws.Range("Y13").Value = 0.8500000000000001
$ws.Range("Z13").Value = -1
$ws.Range("AA13").Value = 0.7749999999999999
$ws.Range("AB13").Value = -1

# Row 21
$ws.Range("B21").Value = 6102188
$ws.Range("E21").Value = 'KR Reykjavik'
$ws.Range("F21").Value = 'Keflavik'
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 'H'
$ws.Range("J21").Value = 1.571
$ws.Range("K21").Value = 4
$ws.Range("L21").Value = 4.75
$ws.Range("M21").Value = 1.65
$ws.Range("N21").Value = 3.75
$ws.Range("O21").Value = 4.5
$ws.Range("P21").Value = -0.75
$ws.Range("Q21").Value = 1.875
$ws.Range("R21").Value = 1.975
$ws.Range("S21").Value = 2.75
$ws.Range("T21").Value = 2
$ws.Range("U21").Value = 1.85
$ws.Range("V21").Value = 0.6499999999999999
$ws.Range("W21").Value = -1
$ws.Range("X21").Value = -1
$ws.Range("Y21").Value = 0.875
$ws.Range("Z21").Value = -1
$ws.Range("AA21").Value = -1
$ws.Range("AB21").Value = 0.8500000000000001

# Row 22
$ws.Range("B22").Value = 6102191
$ws.Range("E22").Value = 'Fram Reykjavik'
$ws.Range("F22").Value = 'HK Kopavogur'
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 2
$ws.Range("I22").Value = 'H'
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = 3.75
$ws.Range("L22").Value = 3
$ws.Range("M22").Value = 1.833
$ws.Range("N22").Value = 3.75
$ws.Range("O22").Value = 3.4
$ws.Range("P22").Value = -0.5
$ws.Range("Q22").Value = 1.875
$ws.Range("R22").Value = 1.975
$ws.Range("S22").Value = 3.25
$ws.Range("T22").Value = 2
$ws.Range("U22").Value = 1.85
$ws.Range("V22").Value = 0.833
$ws.Range("W22").Value = -1
$ws.Range("X22").Value = -1
$ws.Range("Y22").Value = 0.875
$ws.Range("Z22").Value = -1
$ws.Range("AA22").Value = 1
$ws.Range("AB22").Value = -1

# Row 23
$ws.Range("B23").Value = 6102190
$ws.Range("E23").Value = 'Fylkir Reykjavik'
$ws.Range("F23").Value = 'Vikingur Reykjavik'
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 3
$ws.Range("I23").Value = 'A'
$ws.Range("J23").Value = 5.25
$ws.Range("K23").Value = 4.333
$ws.Range("L23").Value = 1.45
$ws.Range("M23").Value = 5.75
$ws.Range("N23").Value = 4.5
$ws.Range("O23").Value = 1.4
$ws.Range("P23").Value = 1.25
$ws.Range("Q23").Value = 1.875
$ws.Range("R23").Value = 1.925
$ws.Range("S23").Value = 3.25
$ws.Range("T23").Value = 1.925
$ws.Range("U23").Value = 1.875
$ws.Range("V23").Value = -1
$ws.Range("W23").Value = -1
$ws.Range("X23").Value = 0.3999999999999999
$ws.Range("Y23").Value = -1
$ws.Range("Z23").Value = 0.925
$ws.Range("AA23").Value = 0.925
$ws.Range("AB23").Value = -1

# Row 24
$ws.Range("B24").Value = 6102189
$ws.Range("E24").Value = 'Stjarnan'
$ws.Range("F24").Value = 'FH Hafnarfjordur'
$ws.Range("G24").Value = 5
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 'H'
$ws.Range("J24").Value = 2.1
$ws.Range("K24").Value = 3.4
$ws.Range("L24").Value = 3
$ws.Range("M24").Value = 2.1
$ws.Range("N24").Value = 3.4
$ws.Range("O24").Value = 2.9
$ws.Range("P24").Value = -0.25
$ws.Range("Q24").Value = 1.925
$ws.Range("R24").Value = 1.875
$ws.Range("S24").Value = 3
$ws.Range("T24").Value = 1.925
$ws.Range("U24").Value = 1.875
$ws.Range("V24").Value = 1.1
$ws.Range("W24").Value = -1
$ws.Range("X24").Value = -1
$ws.Range("Y24").Value = 0.925
$ws.Range("Z24").Value = -1
$ws.Range("AA24").Value = 0.925
$ws.Range("AB24").Value = -1

# Row 33
$ws.Range("B33").Value = 6102200
$ws.Range("E33").Value = 'KR Reykjavik'
$ws.Range("F33").Value = 'FH Hafnarfjordur'
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 'H'
$ws.Range("J33").Value = 1.909
$ws.Range("K33").Value = 3.6
$ws.Range("L33").Value = 3.4
$ws.Range("M33").Value = 2.05
$ws.Range("N33").Value = 3.5
$ws.Range("O33").Value = 3.25
$ws.Range("P33").Value = -0.25
$ws.Range("Q33").Value = 1.8
$ws.Range("R33").Value = 2.05
$ws.Range("S33").Value = 2.75
$ws.Range("T33").Value = 1.925
$ws.Range("U33").Value = 1.925
$ws.Range("V33").Value = 1.05
$ws.Range("W33").Value = -1
$ws.Range("X33").Value = -1
$ws.Range("Y33").Value = 0.8
$ws.Range("Z33").Value = -1
$ws.Range("AA33").Value = -1
$ws.Range("AB33").Value = 0.925

# Row 34
$ws.Range("B34").Value = 6922578
$ws.Range("E34").Value = 'Fylkir Reykjavik'
$ws.Range("F34").Value = 'HK Kopavogur'
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 'D'
$ws.Range("J34").Value = 2.1
$ws.Range("K34").Value = 3.6
$ws.Range("L34").Value = 2.9
$ws.Range("M34").Value = 2.05
$ws.Range("N34").Value = 3.8
$ws.Range("O34").Value = 3
$ws.Range("P34").Value = -0.25
$ws.Range("Q34").Value = 1.85
$ws.Range("R34").Value = 1.95
$ws.Range("S34").Value = 3.5
$ws.Range("T34").Value = 1.925
$ws.Range("U34").Value = 1.775
$ws.Range("V34").Value = -1
$ws.Range("W34").Value = 2.8
$ws.Range("X34").Value = -1
$ws.Range("Y34").Value = -0.5
$ws.Range("Z34").Value = 0.475
$ws.Range("AA34").Value = -1
$ws.Range("AB34").Value = 0.7749999999999999

# Row 45
$ws.Range("B45").Value = 6102350
$ws.Range("E45").Value = 'Keflavik'
$ws.Range("F45").Value = 'FH Hafnarfjordur'
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 3
$ws.Range("I45").Value = 'A'
$ws.Range("J45").Value = 3.8
$ws.Range("K45").Value = 3.6
$ws.Range("L45").Value = 1.909
$ws.Range("M45").Value = 3.8
$ws.Range("N45").Value = 3.75
$ws.Range("O45").Value = 1.85
$ws.Range("P45").Value = 0.5
$ws.Range("Q45").Value = 1.95
$ws.Range("R45").Value = 1.85
$ws.Range("S45").Value = 3
$ws.Range("T45").Value = 1.95
$ws.Range("U45").Value = 1.85
$ws.Range("V45").Value = -1
$ws.Range("W45").Value = -1
$ws.Range("X45").Value = 0.8500000000000001
$ws.Range("Y45").Value = -1
$ws.Range("Z45").Value = 0.8500000000000001
$ws.Range("AA45").Value = 0.95
$ws.Range("AB45").Value = -1

# Row 46
$ws.Range("B46").Value = 6102213
$ws.Range("E46").Value = 'KR Reykjavik'
$ws.Range("F46").Value = 'Valur Reykjavik'
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 4
$ws.Range("I46").Value = 'A'
$ws.Range("J46").Value = 2.875
$ws.Range("K46").Value = 3.5
$ws.Range("L46").Value = 2.3
$ws.Range("M46").Value = 2.7
$ws.Range("N46").Value = 3.6
$ws.Range("O46").Value = 2.45
$ws.Range("P46").Value = 0
$ws.Range("Q46").Value = 2
$ws.Range("R46").Value = 1.85
$ws.Range("S46").Value = 3
$ws.Range("T46").Value = 1.875
$ws.Range("U46").Value = 1.975
$ws.Range("V46").Value = -1
$ws.Range("W46").Value = -1
$ws.Range("X46").Value = 1.45
$ws.Range("Y46").Value = -1
$ws.Range("Z46").Value = 0.8500000000000001
$ws.Range("AA46").Value = 0.875
$ws.Range("AB46").Value = -1

# Row 67
$ws.Range("B67").Value = 6102229
$ws.Range("E67").Value = 'KR Reykjavik'
$ws.Range("F67").Value = 'Fylkir Reykjavik'
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 'H'
$ws.Range("J67").Value = 1.55
$ws.Range("K67").Value = 4.5
$ws.Range("L67").Value = 4.333
$ws.Range("M67").Value = 1.75
$ws.Range("N67").Value = 4.2
$ws.Range("O67").Value = 3.5
$ws.Range("P67").Value = -0.75
$ws.Range("Q67").Value = 1.975
$ws.Range("R67").Value = 1.825
$ws.Range("S67").Value = 3.25
$ws.Range("T67").Value = 1.975
$ws.Range("U67").Value = 1.825
$ws.Range("V67").Value = 0.75
$ws.Range("W67").Value = -1
$ws.Range("X67").Value = -1
$ws.Range("Y67").Value = 0.9750000000000001
$ws.Range("Z67").Value = -1
$ws.Range("AA67").Value = -1
$ws.Range("AB67").Value = 0.825

# Row 68
$ws.Range("B68").Value = 6102230
$ws.Range("E68").Value = 'Keflavik'
$ws.Range("F68").Value = 'Fram Reykjavik'
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 'D'
$ws.Range("J68").Value = 2.55
$ws.Range("K68").Value = 3.4
$ws.Range("L68").Value = 2.45
$ws.Range("M68").Value = 2.55
$ws.Range("N68").Value = 3.3
$ws.Range("O68").Value = 2.45
$ws.Range("P68").Value = 0
$ws.Range("Q68").Value = 2.025
$ws.Range("R68").Value = 1.825
$ws.Range("S68").Value = 3
$ws.Range("T68").Value = 1.9
$ws.Range("U68").Value = 1.95
$ws.Range("V68").Value = -1
$ws.Range("W68").Value = 2.3
$ws.Range("X68").Value = -1
$ws.Range("Y68").Value = 0
$ws.Range("Z68").Value = 0
$ws.Range("AA68").Value = -1
$ws.Range("AB68").Value = 0.95

# Row 72
$ws.Range("B72").Value = 6102357
$ws.Range("E72").Value = 'Valur Reykjavik'
$ws.Range("F72").Value = 'HK Kopavogur'
$ws.Range("G72").Value = 4
$ws.Range("H72").Value = 1
$ws.Range("I72").Value = 'H'
$ws.Range("J72").Value = 1.333
$ws.Range("K72").Value = 5
$ws.Range("L72").Value = 7
$ws.Range("M72").Value = 1.3
$ws.Range("N72").Value = 5.25
$ws.Range("O72").Value = 7
$ws.Range("P72").Value = -1.75
$ws.Range("Q72").Value = 1.975
$ws.Range("R72").Value = 1.825
$ws.Range("S72").Value = 3.75
$ws.Range("T72").Value = 1.875
$ws.Range("U72").Value = 1.925
$ws.Range("V72").Value = 0.3
$ws.Range("W72").Value = -1
$ws.Range("X72").Value = -1
$ws.Range("Y72").Value = 0.9750000000000001
$ws.Range("Z72").Value = -1
$ws.Range("AA72").Value = 0.875
$ws.Range("AB72").Value = -1

# Row 73
$ws.Range("B73").Value = 6102236
$ws.Range("E73").Value = 'IBV Vestmannaeyjar'
$ws.Range("F73").Value = 'KR Reykjavik'
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 2
$ws.Range("I73").Value = 'D'
$ws.Range("J73").Value = 2.625
$ws.Range("K73").Value = 3.6
$ws.Range("L73").Value = 2.25
$ws.Range("M73").Value = 2.625
$ws.Range("N73").Value = 3.6
$ws.Range("O73").Value = 2.25
$ws.Range("P73").Value = 0.25
$ws.Range("Q73").Value = 1.775
$ws.Range("R73").Value = 2.025
$ws.Range("S73").Value = 2.75
$ws.Range("T73").Value = 1.825
$ws.Range("U73").Value = 1.975
$ws.Range("V73").Value = -1
$ws.Range("W73").Value = 2.6
$ws.Range("X73").Value = -1
$ws.Range("Y73").Value = 0.3875
$ws.Range("Z73").Value = -0.5
$ws.Range("AA73").Value = 0.825
$ws.Range("AB73").Value = -1

# Row 74
$ws.Range("B74").Value = 6102235
$ws.Range("E74").Value = 'Breidablik'
$ws.Range("F74").Value = 'FH Hafnarfjordur'
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 2
$ws.Range("I74").Value = 'A'
$ws.Range("J74").Value = 1.6
$ws.Range("K74").Value = 4.2
$ws.Range("L74").Value = 4.333
$ws.Range("M74").Value = 1.65
$ws.Range("N74").Value = 4.2
$ws.Range("O74").Value = 4
$ws.Range("P74").Value = -0.75
$ws.Range("Q74").Value = 1.825
$ws.Range("R74").Value = 1.975
$ws.Range("S74").Value = 3.5
$ws.Range("T74").Value = 1.825
$ws.Range("U74").Value = 1.975
$ws.Range("V74").Value = -1
$ws.Range("W74").Value = -1
$ws.Range("X74").Value = 3
$ws.Range("Y74").Value = -1
$ws.Range("Z74").Value = 0.9750000000000001
$ws.Range("AA74").Value = -1
$ws.Range("AB74").Value = 0.9750000000000001

# Row 75
$ws.Range("B75").Value = 6102234
$ws.Range("E75").Value = 'Fram Reykjavik'
$ws.Range("F75").Value = 'Vikingur Reykjavik'
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 3
$ws.Range("I75").Value = 'A'
$ws.Range("J75").Value = 6
$ws.Range("K75").Value = 5
$ws.Range("L75").Value = 1.4
$ws.Range("M75").Value = 7
$ws.Range("N75").Value = 5
$ws.Range("O75").Value = 1.363
$ws.Range("P75").Value = 1.5
$ws.Range("Q75").Value = 1.95
$ws.Range("R75").Value = 1.9
$ws.Range("S75").Value = 3.5
$ws.Range("T75").Value = 1.9
$ws.Range("U75").Value = 1.95
$ws.Range("V75").Value = -1
$ws.Range("W75").Value = -1
$ws.Range("X75").Value = 0.363
$ws.Range("Y75").Value = 0.95
$ws.Range("Z75").Value = -1
$ws.Range("AA75").Value = 0.8999999999999999
$ws.Range("AB75").Value = -1

# Row 76
$ws.Range("B76").Value = 6102233
$ws.Range("E76").Value = 'Stjarnan'
$ws.Range("F76").Value = 'Keflavik'
$ws.Range("G76").Value = 3
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 'H'
$ws.Range("J76").Value = 1.333
$ws.Range("K76").Value = 5
$ws.Range("L76").Value = 7
$ws.Range("M76").Value = 1.333
$ws.Range("N76").Value = 5
$ws.Range("O76").Value = 7
$ws.Range("P76").Value = -1.5
$ws.Range("Q76").Value = 1.925
$ws.Range("R76").Value = 1.925
$ws.Range("S76").Value = 3.25
$ws.Range("T76").Value = 1.825
$ws.Range("U76").Value = 2.025
$ws.Range("V76").Value = 0.333
$ws.Range("W76").Value = -1
$ws.Range("X76").Value = -1
$ws.Range("Y76").Value = 0.925
$ws.Range("Z76").Value = -1
$ws.Range("AA76").Value = -0.5
$ws.Range("AB76").Value = 0.5125

# Row 85
$ws.Range("B85").Value = 7173163
$ws.Range("E85").Value = 'FH Hafnarfjordur'
$ws.Range("F85").Value = 'Stjarnan'
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 3
$ws.Range("I85").Value = 'A'
$ws.Range("J85").Value = 2.55
$ws.Range("K85").Value = 3.6
$ws.Range("L85").Value = 2.375
$ws.Range("M85").Value = 2.7
$ws.Range("N85").Value = 3.5
$ws.Range("O85").Value = 2.25
$ws.Range("P85").Value = 0.25
$ws.Range("Q85").Value = 1.775
$ws.Range("R85").Value = 2.025
$ws.Range("S85").Value = 3
$ws.Range("T85").Value = 1.825
$ws.Range("U85").Value = 1.975
$ws.Range("V85").Value = -1
$ws.Range("W85").Value = -1
$ws.Range("X85").Value = 1.25
$ws.Range("Y85").Value = -1
$ws.Range("Z85").Value = 1.025
$ws.Range("AA85").Value = 0.825
$ws.Range("AB85").Value = -1

# Row 86
$ws.Range("B86").Value = 7173164
$ws.Range("E86").Value = 'KR Reykjavik'
$ws.Range("F86").Value = 'Valur Reykjavik'
$ws.Range("G86").Value = 2
$ws.Range("H86").Value = 2
$ws.Range("I86").Value = 'D'
$ws.Range("J86").Value = 2.5
$ws.Range("K86").Value = 3.6
$ws.Range("L86").Value = 2.375
$ws.Range("M86").Value = 2.25
$ws.Range("N86").Value = 3.8
$ws.Range("O86").Value = 2.6
$ws.Range("P86").Value = -0.25
$ws.Range("Q86").Value = 1.95
$ws.Range("R86").Value = 1.75
$ws.Range("S86").Value = 3.5
$ws.Range("T86").Value = 1.825
$ws.Range("U86").Value = 1.975
$ws.Range("V86").Value = -1
$ws.Range("W86").Value = 2.8
$ws.Range("X86").Value = -1
$ws.Range("Y86").Value = -0.5
$ws.Range("Z86").Value = 0.375
$ws.Range("AA86").Value = 0.825
$ws.Range("AB86").Value = -1

# Row 91
$ws.Range("B91").Value = 7173183
$ws.Range("E91").Value = 'HK Kopavogur'
$ws.Range("F91").Value = 'Fylkir Reykjavik'
$ws.Range("G91").Value = 2
$ws.Range("H91").Value = 2
$ws.Range("I91").Value = 'D'
$ws.Range("J91").Value = 2.2
$ws.Range("K91").Value = 3.6
$ws.Range("L91").Value = 2.7
$ws.Range("M91").Value = 2.45
$ws.Range("N91").Value = 3.6
$ws.Range("O91").Value = 2.4
$ws.Range("P91").Value = 0
$ws.Range("Q91").Value = 1.925
$ws.Range("R91").Value = 1.925
$ws.Range("S91").Value = 3.25
$ws.Range("T91").Value = 1.9
$ws.Range("U91").Value = 1.95
$ws.Range("V91").Value = -1
$ws.Range("W91").Value = 2.6
$ws.Range("X91").Value = -1
$ws.Range("Y91").Value = 0
$ws.Range("Z91").Value = 0
$ws.Range("AA91").Value = 0.8999999999999999
$ws.Range("AB91").Value = -1

# Row 92
$ws.Range("B92").Value = 7173168
$ws.Range("E92").Value = 'Stjarnan'
$ws.Range("F92").Value = 'KR Reykjavik'
$ws.Range("G92").Value = 2
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 'H'
$ws.Range("J92").Value = 1.95
$ws.Range("K92").Value = 3.75
$ws.Range("L92").Value = 3.2
$ws.Range("M92").Value = 1.85
$ws.Range("N92").Value = 3.8
$ws.Range("O92").Value = 3.4
$ws.Range("P92").Value = -0.5
$ws.Range("Q92").Value = 1.825
$ws.Range("R92").Value = 1.975
$ws.Range("S92").Value = 3.25
$ws.Range("T92").Value = 1.825
$ws.Range("U92").Value = 1.975
$ws.Range("V92").Value = 0.8500000000000001
$ws.Range("W92").Value = -1
$ws.Range("X92").Value = -1
$ws.Range("Y92").Value = 0.825
$ws.Range("Z92").Value = -1
$ws.Range("AA92").Value = -1
$ws.Range("AB92").Value = 0.9750000000000001

# Row 93
$ws.Range("B93").Value = 7173166
$ws.Range("E93").Value = 'Valur Reykjavik'
$ws.Range("F93").Value = 'Breidablik'
$ws.Range("G93").Value = 4
$ws.Range("H93").Value = 2
$ws.Range("I93").Value = 'H'
$ws.Range("J93").Value = 2.25
$ws.Range("K93").Value = 3.6
$ws.Range("L93").Value = 2.6
$ws.Range("M93").Value = 2.7
$ws.Range("N93").Value = 4
$ws.Range("O93").Value = 2.1
$ws.Range("P93").Value = 0.25
$ws.Range("Q93").Value = 1.925
$ws.Range("R93").Value = 1.875
$ws.Range("S93").Value = 3.75
$ws.Range("T93").Value = 1.925
$ws.Range("U93").Value = 1.775
$ws.Range("V93").Value = 1.7
$ws.Range("W93").Value = -1
$ws.Range("X93").Value = -1
$ws.Range("Y93").Value = 0.925
$ws.Range("Z93").Value = -1
$ws.Range("AA93").Value = 0.925
$ws.Range("AB93").Value = -1

# Row 94
$ws.Range("B94").Value = 7173182
$ws.Range("E94").Value = 'Fram Reykjavik'
$ws.Range("F94").Value = 'Keflavik'
$ws.Range("G94").Value = 3
$ws.Range("H94").Value = 1
$ws.Range("I94").Value = 'H'
$ws.Range("J94").Value = 1.75
$ws.Range("K94").Value = 3.8
$ws.Range("L94").Value = 3.75
$ws.Range("M94").Value = 1.833
$ws.Range("N94").Value = 3.8
$ws.Range("O94").Value = 3.5
$ws.Range("P94").Value = -0.5
$ws.Range("Q94").Value = 1.825
$ws.Range("R94").Value = 1.975
$ws.Range("S94").Value = 3.25
$ws.Range("T94").Value = 1.925
$ws.Range("U94").Value = 1.875
$ws.Range("V94").Value = 0.833
$ws.Range("W94").Value = -1
$ws.Range("X94").Value = -1
$ws.Range("Y94").Value = 0.825
$ws.Range("Z94").Value = -1
$ws.Range("AA94").Value = 0.925
$ws.Range("AB94").Value = -1

# Row 95
$ws.Range("B95").Value = 7173167
$ws.Range("E95").Value = 'Vikingur Reykjavik'
$ws.Range("F95").Value = 'FH Hafnarfjordur'
$ws.Range("G95").Value = 2
$ws.Range("H95").Value = 1
$ws.Range("I95").Value = 'H'
$ws.Range("J95").Value = 1.6
$ws.Range("K95").Value = 4
$ws.Range("L95").Value = 4.333
$ws.Range("M95").Value = 1.571
$ws.Range("N95").Value = 4.2
$ws.Range("O95").Value = 4.5
$ws.Range("P95").Value = -1
$ws.Range("Q95").Value = 1.875
$ws.Range("R95").Value = 1.925
$ws.Range("S95").Value = 3.25
$ws.Range("T95").Value = 1.775
$ws.Range("U95").Value = 1.925
$ws.Range("V95").Value = 0.571
$ws.Range("W95").Value = -1
$ws.Range("X95").Value = -1
$ws.Range("Y95").Value = 0
$ws.Range("Z95").Value = 0
$ws.Range("AA95").Value = -0.5
$ws.Range("AB95").Value = 0.4625

# Row 98
$ws.Range("B98").Value = 7173185
$ws.Range("E98").Value = 'HK Kopavogur'
$ws.Range("F98").Value = 'IBV Vestmannaeyjar'
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 1
$ws.Range("I98").Value = 'A'
$ws.Range("J98").Value = 2.625
$ws.Range("K98").Value = 3.6
$ws.Range("L98").Value = 2.3
$ws.Range("M98").Value = 2.7
$ws.Range("N98").Value = 3.6
$ws.Range("O98").Value = 2.25
$ws.Range("P98").Value = 0.25
$ws.Range("Q98").Value = 1.825
$ws.Range("R98").Value = 1.975
$ws.Range("S98").Value = 3.25
$ws.Range("T98").Value = 1.925
$ws.Range("U98").Value = 1.875
$ws.Range("V98").Value = -1
$ws.Range("W98").Value = -1
$ws.Range("X98").Value = 1.25
$ws.Range("Y98").Value = -1
$ws.Range("Z98").Value = 0.9750000000000001
$ws.Range("AA98").Value = -1
$ws.Range("AB98").Value = 0.875

# Row 99
$ws.Range("B99").Value = 7173186
$ws.Range("E99").Value = 'Fram Reykjavik'
$ws.Range("F99").Value = 'KA Akureyri'
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 'H'
$ws.Range("J99").Value = 2.5
$ws.Range("K99").Value = 3.75
$ws.Range("L99").Value = 2.375
$ws.Range("M99").Value = 2
$ws.Range("N99").Value = 3.8
$ws.Range("O99").Value = 3
$ws.Range("P99").Value = -0.25
$ws.Range("Q99").Value = 1.825
$ws.Range("R99").Value = 2.025
$ws.Range("S99").Value = 3.5
$ws.Range("T99").Value = 1.925
$ws.Range("U99").Value = 1.925
$ws.Range("V99").Value = 1
$ws.Range("W99").Value = -1
$ws.Range("X99").Value = -1
$ws.Range("Y99").Value = 0.825
$ws.Range("Z99").Value = -1
$ws.Range("AA99").Value = -1
$ws.Range("AB99").Value = 0.925

# Row 103
$ws.Range("B103").Value = 7173172
$ws.Range("E103").Value = 'FH Hafnarfjordur'
$ws.Range("F103").Value = 'KR Reykjavik'
$ws.Range("G103").Value = 3
$ws.Range("H103").Value = 1
$ws.Range("I103").Value = 'H'
$ws.Range("J103").Value = 2.2
$ws.Range("K103").Value = 3.6
$ws.Range("L103").Value = 2.6
$ws.Range("M103").Value = 2.25
$ws.Range("N103").Value = 3.75
$ws.Range("O103").Value = 2.55
$ws.Range("P103").Value = 0
$ws.Range("Q103").Value = 1.775
$ws.Range("R103").Value = 2.025
$ws.Range("S103").Value = 3.5
$ws.Range("T103").Value = 1.875
$ws.Range("U103").Value = 1.925
$ws.Range("V103").Value = 1.25
$ws.Range("W103").Value = -1
$ws.Range("X103").Value = -1
$ws.Range("Y103").Value = 0.7749999999999999
$ws.Range("Z103").Value = -1
$ws.Range("AA103").Value = 0.875
$ws.Range("AB103").Value = -1

# Row 104
$ws.Range("B104").Value = 7173188
$ws.Range("E104").Value = 'IBV Vestmannaeyjar'
$ws.Range("F104").Value = 'Keflavik'
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 1
$ws.Range("I104").Value = 'D'
$ws.Range("J104").Value = 1.4
$ws.Range("K104").Value = 4.5
$ws.Range("L104").Value = 5.5
$ws.Range("M104").Value = 1.222
$ws.Range("N104").Value = 5.5
$ws.Range("O104").Value = 9
$ws.Range("P104").Value = -2
$ws.Range("Q104").Value = 1.975
$ws.Range("R104").Value = 1.825
$ws.Range("S104").Value = 3.75
$ws.Range("T104").Value = 1.9
$ws.Range("U104").Value = 1.9
$ws.Range("V104").Value = -1
$ws.Range("W104").Value = 4.5
$ws.Range("X104").Value = -1
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = 0.825
$ws.Range("AA104").Value = -1
$ws.Range("AB104").Value = 0.8999999999999999

# Row 106
$ws.Range("B106").Value = 7173187
$ws.Range("E106").Value = 'Fylkir Reykjavik'
$ws.Range("F106").Value = 'Fram Reykjavik'
$ws.Range("G106").Value = 5
$ws.Range("H106").Value = 1
$ws.Range("I106").Value = 'H'
$ws.Range("J106").Value = 2
$ws.Range("K106").Value = 3.75
$ws.Range("L106").Value = 2.9
$ws.Range("M106").Value = 2.1
$ws.Range("N106").Value = 3.5
$ws.Range("O106").Value = 3
$ws.Range("P106").Value = -0.25
$ws.Range("Q106").Value = 1.825
$ws.Range("R106").Value = 2.025
$ws.Range("S106").Value = 3.5
$ws.Range("T106").Value = 2
$ws.Range("U106").Value = 1.85
$ws.Range("V106").Value = 1.1
$ws.Range("W106").Value = -1
$ws.Range("X106").Value = -1
$ws.Range("Y106").Value = 0.825
$ws.Range("Z106").Value = -1
$ws.Range("AA106").Value = 1
$ws.Range("AB106").Value = -1
